$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to be treated as text so that
# numeric-looking values (e.g. "1.00", "0.710") keep their exact
# original formatting instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range("D2").Value = '70.851.55'
$ws.Range("D3").Value = '3.782.76'
$ws.Range("D5").Value = '697.70'
$ws.Range("D6").Value = '169.00'
$ws.Range("D7").Value = '3.783.43'
$ws.Range("D8").Value = '1.00'
$ws.Range("D14").Value = '35.92'
$ws.Range("D15").Value = '4.425.26'
$ws.Range("D16").Value = '3.826.49'
$ws.Range("D17").Value = '71.042.52'
$ws.Range("D19").Value = '17.45'
$ws.Range("D20").Value = '7.13'
$ws.Range("D21").Value = '515.06'
$ws.Range("D22").Value = '10.33'
$ws.Range("D23").Value = '0.710'
$ws.Range("D26").Value = '12.51'
$ws.Range("D27").Value = '3.936.79'
$ws.Range("D28").Value = '10.16'
$ws.Range("D30").Value = '1.94'
$ws.Range("D34").Value = '28.95'
$ws.Range("D36").Value = '9.19'
$ws.Range("D38").Value = '3.747.52'
$ws.Range("D40").Value = '0.0994'
$ws.Range("D41").Value = '2.34'
$ws.Range("D46").Value = '163.96'
$ws.Range("D47").Value = '49.12'
$ws.Range("D49").Value = '416.13'
$ws.Range("D50").Value = '8.58'
$ws.Range("D51").Value = '1.36'

# Restore the original (default) style on column D so no stray
# number-format is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("E3").Value = '  -1.29%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("E6").Value = '  -1.53%  '
$ws.Range("E7").Value = '  -1.15%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E10").Value = '  -1.96%  '
$ws.Range("E11").Value = '  +2.78%  '
$ws.Range("E12").Value = '  +3.69%  '
$ws.Range("E13").Value = '  -2.48%  '
$ws.Range("E14").Value = '  -2.11%  '
$ws.Range("E15").Value = '  -1.21%  '
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("E20").Value = '  -1.35%  '
$ws.Range("E21").Value = '  +3.66%  '
$ws.Range("E22").Value = '  -3.37%  '
$ws.Range("E23").Value = '  -3.39%  '
$ws.Range("E24").Value = '  -2.00%  '
$ws.Range("E25").Value = '  -3.56%  '
$ws.Range("E26").Value = '  +2.94%  '
$ws.Range("E27").Value = '  -1.20%  '
$ws.Range("E28").Value = '  -3.97%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  -7.07%  '
$ws.Range("E31").Value = '  -4.35%  '
$ws.Range("E32").Value = '  +0.25%  '
$ws.Range("E33").Value = '  -2.46%  '
$ws.Range("E34").Value = '  -1.58%  '
$ws.Range("E35").Value = '  -3.68%  '
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("E38").Value = '  -1.24%  '
$ws.Range("E39").Value = '  +9.71%  '
$ws.Range("E40").Value = '  -2.94%  '
$ws.Range("E41").Value = '  +1.13%  '
$ws.Range("E42").Value = '  -2.81%  '
$ws.Range("E44").Value = '  -4.47%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("E48").Value = '  -4.72%  '
$ws.Range("E49").Value = '  -3.55%  '
$ws.Range("E50").Value = '  -1.80%  '
$ws.Range("E51").Value = '  -1.19%  '

